$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2986
$ws.Range("I70").Value = 3325
$ws.Range("J70").Value = 2910.6667
$ws.Range("K70").Value = 9975
$ws.Range("L70").Value = 8732.000100000001
$ws.Range("M70").Value = -9705
$ws.Range("N70").Value = -9272.000100000001

# Row 73
$ws.Range("H73").Value = 2986
$ws.Range("I73").Value = 3325
$ws.Range("J73").Value = 2910.6667
$ws.Range("K73").Value = 9975
$ws.Range("L73").Value = 8732.000100000001
$ws.Range("M73").Value = -9039
$ws.Range("N73").Value = -10604.0001

# Row 98
$ws.Range("H98").Value = 1125.4517
$ws.Range("I98").Value = 1205.619
$ws.Range("J98").Value = 957.1
$ws.Range("K98").Value = 1205.619
$ws.Range("L98").Value = 957.1
$ws.Range("M98").Value = 292.3810000000001
$ws.Range("N98").Value = -3953.1

# Row 122
$ws.Range("H122").Value = 1125.4517
$ws.Range("I122").Value = 1205.619
$ws.Range("J122").Value = 957.1
$ws.Range("K122").Value = 3616.857
$ws.Range("L122").Value = 2871.3
$ws.Range("M122").Value = -1166.857
$ws.Range("N122").Value = -7771.3

# Row 135
$ws.Range("H135").Value = 15626840
$ws.Range("I135").Value = 2314.5
$ws.Range("J135").Value = 22728898
$ws.Range("K135").Value = 20830.5
$ws.Range("L135").Value = 204560082
$ws.Range("M135").Value = -18295.5
$ws.Range("N135").Value = -204565152

# Row 137
$ws.Range("H137").Value = 4249.4116
$ws.Range("I137").Value = 1053.2858
$ws.Range("J137").Value = 5078.037
$ws.Range("K137").Value = 3159.8574
$ws.Range("L137").Value = 15234.111
$ws.Range("M137").Value = -609.8574000000003
$ws.Range("N137").Value = -20334.111

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4122.3
$ws.Range("I61").Value = 2469.2222
$ws.Range("J61").Value = 19000
$ws.Range("K61").Value = 2469.2222
$ws.Range("L61").Value = 19000
$ws.Range("M61").Value = -2257.2222
$ws.Range("N61").Value = -19424

# Row 74
$ws.Range("H74").Value = 66165.516
$ws.Range("I74").Value = 81613.24000000001
$ws.Range("K74").Value = 81613.24000000001
$ws.Range("M74").Value = -80739.24000000001

# Row 77
$ws.Range("H77").Value = 66165.516
$ws.Range("I77").Value = 81613.24000000001
$ws.Range("K77").Value = 408066.2
$ws.Range("M77").Value = -403698.2

# Row 102
$ws.Range("H102").Value = 2833
$ws.Range("I102").Value = 2833
$ws.Range("K102").Value = 2833
$ws.Range("M102").Value = -1211

# Row 132
$ws.Range("H132").Value = 2138193.5
$ws.Range("I132").Value = 2406122.2
$ws.Range("J132").Value = 1012892
$ws.Range("K132").Value = 7218366.600000001
$ws.Range("L132").Value = 3038676
$ws.Range("M132").Value = -7215836.600000001
$ws.Range("N132").Value = -3043736

# Row 136
$ws.Range("H136").Value = 4122.3
$ws.Range("I136").Value = 2469.2222
$ws.Range("J136").Value = 19000
$ws.Range("K136").Value = 7407.6666
$ws.Range("L136").Value = 57000
$ws.Range("M136").Value = -4857.6666
$ws.Range("N136").Value = -62100

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 28481.29
$ws.Range("I134").Value = 1527.6786
$ws.Range("J134").Value = 103951.4
$ws.Range("K134").Value = 4583.0358
$ws.Range("L134").Value = 311854.2
$ws.Range("M134").Value = -2048.0358
$ws.Range("N134").Value = -316924.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2467.6667
$ws.Range("I31").Value = 988.8929000000001
$ws.Range("K31").Value = 988.8929000000001
$ws.Range("M31").Value = -693.8929000000001

# Row 34
$ws.Range("H34").Value = 2467.6667
$ws.Range("I34").Value = 988.8929000000001
$ws.Range("K34").Value = 988.8929000000001
$ws.Range("M34").Value = -786.8929000000001

# Row 58
$ws.Range("H58").Value = 3201.037
$ws.Range("I58").Value = 1009.7143
$ws.Range("J58").Value = 5560.923
$ws.Range("K58").Value = 1009.7143
$ws.Range("L58").Value = 5560.923
$ws.Range("M58").Value = -806.7143
$ws.Range("N58").Value = -5966.923

# Row 132
$ws.Range("H132").Value = 2069.2307
$ws.Range("I132").Value = 1340.4706
$ws.Range("K132").Value = 4021.4118
$ws.Range("M132").Value = -1491.4118

# Row 134
$ws.Range("H134").Value = 12822061
$ws.Range("I134").Value = 1303.5186
$ws.Range("J134").Value = 41668764
$ws.Range("K134").Value = 3910.5558
$ws.Range("L134").Value = 125006292
$ws.Range("M134").Value = -1375.5558
$ws.Range("N134").Value = -125011362

# Row 136
$ws.Range("H136").Value = 3201.037
$ws.Range("I136").Value = 1009.7143
$ws.Range("J136").Value = 5560.923
$ws.Range("K136").Value = 3029.1429
$ws.Range("L136").Value = 16682.769
$ws.Range("M136").Value = -479.1428999999998
$ws.Range("N136").Value = -21782.769

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 200.1875
$ws.Range("I40").Value = 129.41667
$ws.Range("J40").Value = 412.5
$ws.Range("K40").Value = 517.66668
$ws.Range("L40").Value = 1650
$ws.Range("M40").Value = -448.66668
$ws.Range("N40").Value = -1788

$ws = $wb.Worksheets.Item("GSM")
# Row 127
$ws.Range("H127").Value = 40678.25
$ws.Range("J127").Value = 40678.25
$ws.Range("L127").Value = 40678.25
$ws.Range("N127").Value = -50598.25

# Row 132
$ws.Range("H132").Value = 2384.4614
$ws.Range("I132").Value = 1320.1482
$ws.Range("J132").Value = 4779.1665
$ws.Range("K132").Value = 3960.4446
$ws.Range("L132").Value = 14337.4995
$ws.Range("M132").Value = -1430.4446
$ws.Range("N132").Value = -19397.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3074.4666
$ws.Range("I7").Value = 3711
$ws.Range("J7").Value = 2347
$ws.Range("K7").Value = 3711
$ws.Range("L7").Value = 2347
$ws.Range("M7").Value = -3599
$ws.Range("N7").Value = -2571

# Row 126
$ws.Range("H126").Value = 3074.4666
$ws.Range("I126").Value = 3711
$ws.Range("J126").Value = 2347
$ws.Range("K126").Value = 11133
$ws.Range("L126").Value = 7041
$ws.Range("M126").Value = -8663
$ws.Range("N126").Value = -11981

# Row 132
$ws.Range("H132").Value = 347603.3
$ws.Range("I132").Value = 2658.8235
$ws.Range("K132").Value = 7976.470499999999
$ws.Range("M132").Value = -5446.470499999999

# Row 136
$ws.Range("H136").Value = 2341.8
$ws.Range("I136").Value = 1137.25
$ws.Range("J136").Value = 3718.4285
$ws.Range("K136").Value = 3411.75
$ws.Range("L136").Value = 11155.2855
$ws.Range("M136").Value = -861.75
$ws.Range("N136").Value = -16255.2855

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1439.4103
$ws.Range("I132").Value = 1127.6552
$ws.Range("K132").Value = 3382.9656
$ws.Range("M132").Value = -852.9655999999995

# Row 136
$ws.Range("H136").Value = 3761464
$ws.Range("I136").Value = 3761464
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11284392
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11281842
$ws.Range("N136").ClearContents()
